$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 49-53: "E1 HET; E2 cKO" -> "E1 HET; E2cKO"
$ws.Range("B49").Value2 = "E1 HET; E2cKO"
$ws.Range("B50").Value2 = "E1 HET; E2cKO"
$ws.Range("B51").Value2 = "E1 HET; E2cKO"
$ws.Range("B52").Value2 = "E1 HET; E2cKO"
$ws.Range("B53").Value2 = "E1 HET; E2cKO"

# Rows 54-59: "E1 ko; E2 HET" -> "E1ko; E2 HET"
$ws.Range("B54").Value2 = "E1ko; E2 HET"
$ws.Range("B55").Value2 = "E1ko; E2 HET"
$ws.Range("B56").Value2 = "E1ko; E2 HET"
$ws.Range("B57").Value2 = "E1ko; E2 HET"
$ws.Range("B58").Value2 = "E1ko; E2 HET"
$ws.Range("B59").Value2 = "E1ko; E2 HET"

# Update the active selection from B62 to B61
$ws.Range("B61").Select()
